$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4. Set cells left-to-right so new shared-string entries are
# appended in the same order the source workbook used.
#
# Cells whose value "looks like a number" (pure digits, possibly with a
# leading zero) get a leading apostrophe so Excel stores them as text
# instead of silently coercing to a Number (which would also lose leading
# zeros / precision on the 18-digit ID). ClearFormats() afterwards removes
# the resulting "quote prefix" cell style so the cell keeps default
# (unstyled) formatting, matching the source data which was authored with
# the cells as plain shared-string text and no explicit style.
$ws.Range("A4").Value = "'410067"
$ws.Range("D4").Value = "试一试"
$ws.Range("E4").Value = "男"
$ws.Range("F4").Value = "'19931214"
$ws.Range("G4").Value = "'1"
$ws.Range("H4").Value = "'410304199312140590"
$ws.Range("I4").Value = "'04"
$ws.Range("J4").Value = "'01"
$ws.Range("K4").Value = "'4"
$ws.Range("L4").Value = "'2"
$ws.Range("M4").Value = "'1"
$ws.Range("R4").Value = "'13298309877"
$ws.Range("T4").Value = "43-2012215-201405598"

$ws.Range("A4").ClearFormats()
$ws.Range("F4").ClearFormats()
$ws.Range("G4").ClearFormats()
$ws.Range("H4").ClearFormats()
$ws.Range("I4").ClearFormats()
$ws.Range("J4").ClearFormats()
$ws.Range("K4").ClearFormats()
$ws.Range("L4").ClearFormats()
$ws.Range("M4").ClearFormats()
$ws.Range("R4").ClearFormats()
